# Update detik.com news index sheet to the latest scrape snapshot:
#  - rows 2-21: refresh headline/time-ago/url for the 20 newest articles
#  - rows 22-143: bump the "tanggal_berita" footer-link date from 2025-09-30 to 2025-10-01
#  - rows 144-163: drop the now-stale tail of yesterday's articles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-21: new headline (A), relative time (B), article url (D) ---
$ws.Cells.Item(2,1).Value = 'Upacara Kesaktian Pancasila: Muzani Bacakan Teks Pancasila, Puan Baca Ikrar'
$ws.Cells.Item(2,2).Value = '1 menit yang lalu'
$ws.Cells.Item(2,4).Value = 'https://news.detik.com/berita/d-8138859/upacara-kesaktian-pancasila-muzani-bacakan-teks-pancasila-puan-baca-ikrar'

$ws.Cells.Item(3,1).Value = 'Heboh Bobby Minta BL Jadi BK, Legislator: Harusnya Persuasif Agar Tak Gaduh'
$ws.Cells.Item(3,2).Value = '6 menit yang lalu'
$ws.Cells.Item(3,4).Value = 'https://news.detik.com/berita/d-8138857/heboh-bobby-minta-bl-jadi-bk-legislator-harusnya-persuasif-agar-tak-gaduh'

$ws.Cells.Item(4,1).Value = 'Lestari Moerdijat Tekankan Pentingnya Keseimbangan Antara Ilmu dan Iman'
$ws.Cells.Item(4,2).Value = '26 menit yang lalu'
$ws.Cells.Item(4,4).Value = 'https://news.detik.com/berita/d-8138841/lestari-moerdijat-tekankan-pentingnya-keseimbangan-antara-ilmu-dan-iman'

$ws.Cells.Item(5,1).Value = 'Dampak Gempa M 6,5 di Sumenep: 22 Rumah Rusak, 3 Orang Luka'
$ws.Cells.Item(5,2).Value = '27 menit yang lalu'
$ws.Cells.Item(5,4).Value = 'https://news.detik.com/berita/d-8138840/dampak-gempa-m-6-5-di-sumenep-22-rumah-rusak-3-orang-luka'

$ws.Cells.Item(6,1).Value = 'Lika-liku Mercy BJ Habibie: Dicicil RK, Disita KPK, Balik ke Anaknya'
$ws.Cells.Item(6,2).Value = '28 menit yang lalu'
$ws.Cells.Item(6,4).Value = 'https://news.detik.com/berita/d-8138837/lika-liku-mercy-bj-habibie-dicicil-rk-disita-kpk-balik-ke-anaknya'

$ws.Cells.Item(7,1).Value = 'Prabowo Perdana Pimpin Upacara Hari Kesaktian Pancasila'
$ws.Cells.Item(7,2).Value = '41 menit yang lalu'
$ws.Cells.Item(7,4).Value = 'https://news.detik.com/berita/d-8138833/prabowo-perdana-pimpin-upacara-hari-kesaktian-pancasila'

$ws.Cells.Item(8,1).Value = '4 Dinding Rumah Warga di Sumenep Roboh Usai Diguncang Gempa M 6,5'
$ws.Cells.Item(8,2).Value = '43 menit yang lalu'
$ws.Cells.Item(8,4).Value = 'https://news.detik.com/berita/d-8138831/4-dinding-rumah-warga-di-sumenep-roboh-usai-diguncang-gempa-m-6-5'

$ws.Cells.Item(9,1).Value = '4 Pernyataan Keluarga Diplomat Arya Daru Blak-blakan di Senayan'
$ws.Cells.Item(9,2).Value = '55 menit yang lalu'
$ws.Cells.Item(9,4).Value = 'https://news.detik.com/berita/d-8138818/4-pernyataan-keluarga-diplomat-arya-daru-blak-blakan-di-senayan'

$ws.Cells.Item(10,1).Value = 'Respons Kepala BGN soal Ortu Murid SDIT Al Izzah Serang Tolak Diberi MBG'
$ws.Cells.Item(10,2).Value = '1 jam yang lalu'
$ws.Cells.Item(10,4).Value = 'https://news.detik.com/berita/d-8138798/respons-kepala-bgn-soal-ortu-murid-sdit-al-izzah-serang-tolak-diberi-mbg'

$ws.Cells.Item(11,1).Value = 'Musala Ponpes Roboh Tewaskan 3 Santri, Anggota DPR: Yang Lalai Harus Disanksi'
$ws.Cells.Item(11,2).Value = '1 jam yang lalu'
$ws.Cells.Item(11,4).Value = 'https://news.detik.com/berita/d-8138794/musala-ponpes-roboh-tewaskan-3-santri-anggota-dpr-yang-lalai-harus-disanksi'

$ws.Cells.Item(12,1).Value = 'Dendam Sering Dibully, Siswa SMP di Lampung Bunuh Teman Sekolah'
$ws.Cells.Item(12,2).Value = '1 jam yang lalu'
$ws.Cells.Item(12,4).Value = 'https://news.detik.com/berita/d-8138791/dendam-sering-dibully-siswa-smp-di-lampung-bunuh-teman-sekolah'

$ws.Cells.Item(13,1).Value = 'Akhirnya Bocah Korban ''Ayah Juna'' Kembali ke Keluarga'
$ws.Cells.Item(13,2).Value = '1 jam yang lalu'
$ws.Cells.Item(13,4).Value = 'https://news.detik.com/berita/d-8138770/akhirnya-bocah-korban-ayah-juna-kembali-ke-keluarga'

$ws.Cells.Item(14,1).Value = 'Komisi X DPR Minta Ada Aturan Cegah Smart TV Bantuan Pusat Dipakai Karaoke'
$ws.Cells.Item(14,2).Value = '2 jam yang lalu'
$ws.Cells.Item(14,4).Value = 'https://news.detik.com/berita/d-8138766/komisi-x-dpr-minta-ada-aturan-cegah-smart-tv-bantuan-pusat-dipakai-karaoke'

$ws.Cells.Item(15,1).Value = 'Hari Lanjut Usia Internasional 2025: Latar Belakang dan Tema'
$ws.Cells.Item(15,2).Value = '2 jam yang lalu'
$ws.Cells.Item(15,4).Value = 'https://news.detik.com/berita/d-8135251/hari-lanjut-usia-internasional-2025-latar-belakang-dan-tema'

$ws.Cells.Item(16,1).Value = 'Drama Razman ke Luar Negeri Tanpa Izin tapi Hakim Tetap Baca Putusan'
$ws.Cells.Item(16,2).Value = '2 jam yang lalu'
$ws.Cells.Item(16,4).Value = 'https://news.detik.com/berita/d-8138746/drama-razman-ke-luar-negeri-tanpa-izin-tapi-hakim-tetap-baca-putusan'

$ws.Cells.Item(17,1).Value = 'Upacara Hari Kesaktian Pancasila, Apakah Ada Pengibaran Bendera?'
$ws.Cells.Item(17,2).Value = '2 jam yang lalu'
$ws.Cells.Item(17,4).Value = 'https://news.detik.com/berita/d-8138179/upacara-hari-kesaktian-pancasila-apakah-ada-pengibaran-bendera'

$ws.Cells.Item(18,1).Value = 'Komisi VII DPR Harap RUU Kepariwisataan Bisa Disahkan Jadi UU Besok'
$ws.Cells.Item(18,2).Value = '2 jam yang lalu'
$ws.Cells.Item(18,4).Value = 'https://news.detik.com/berita/d-8138737/komisi-vii-dpr-harap-ruu-kepariwisataan-bisa-disahkan-jadi-uu-besok'

$ws.Cells.Item(19,1).Value = '1 Oktober 2025 Memperingati Hari Apa? Ini Daftar Hari Pentingnya'
$ws.Cells.Item(19,2).Value = '2 jam yang lalu'
$ws.Cells.Item(19,4).Value = 'https://news.detik.com/berita/d-8137161/1-oktober-2025-memperingati-hari-apa-ini-daftar-hari-pentingnya'

$ws.Cells.Item(20,1).Value = 'Geger Guru Pandeglang Asyik Karaoke Pakai Smart TV Bantuan Pusat'
$ws.Cells.Item(20,2).Value = '2 jam yang lalu'
$ws.Cells.Item(20,4).Value = 'https://news.detik.com/berita/d-8138734/geger-guru-pandeglang-asyik-karaoke-pakai-smart-tv-bantuan-pusat'

$ws.Cells.Item(21,1).Value = 'Satgas Damai Cartenz Bongkar Pemasok Senpi-Amunisi KKB, 2 Orang Diciduk'
$ws.Cells.Item(21,2).Value = '3 jam yang lalu'
$ws.Cells.Item(21,4).Value = 'https://news.detik.com/berita/d-8138728/satgas-damai-cartenz-bongkar-pemasok-senpi-amunisi-kkb-2-orang-diciduk'

# --- Rows 22-143: only column B (tanggal_berita) changes, 2025-09-30 -> 2025-10-01.
# Leading apostrophe keeps it stored as literal text (matches the original "General"
# formatted text cell) instead of Excel auto-coercing the ISO-looking string to a date. ---
for ($r = 22; $r -le 143; $r++) {
    $ws.Cells.Item($r,2).Value = "'2025-10-01"
}

# --- Rows 144-163 (yesterday's oldest articles) are dropped entirely;
# dimension shrinks from A1:F163 to A1:F143. ---
$ws.Range("A144:F163").EntireRow.Delete()
